$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.4551375938158524
$ws.Range("J4").Value = 0.4962935625645516
$ws.Range("K4").Value = 0.4459566511720581
$ws.Range("L4").Value = 2.658316893901198
